# Update the "cryptos" worksheet with refreshed market data.
# Mirrors the GitHub Actions scrape that refreshed Price / Volume(1h)
# figures (and, for rows 40/41, a rank swap between Aave and Kaspa).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '76.323.55'

# Row 3 - Ethereum
$ws.Range("D3").Value = '3.039.83'
$ws.Range("E3").Value = '  +3.79%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  -0.03%  '

# Row 5 - Solana
$ws.Range("D5").Value = "'200.16"
$ws.Range("E5").Value = '  -1.44%  '

# Row 6 - BNB
$ws.Range("D6").Value = "'623.86"
$ws.Range("E6").Value = '  +4.60%  '

# Row 7 - USDC
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = '  +0.04%  '

# Row 8 - XRP
$ws.Range("E8").Value = '  -0.50%  '

# Row 9 - Dogecoin
$ws.Range("E9").Value = '  +4.72%  '

# Row 10 - LidoStakedEther
$ws.Range("D10").Value = '3.040.31'
$ws.Range("E10").Value = '  +3.86%  '

# Row 11 - Cardano
$ws.Range("D11").Value = "'0.440"
$ws.Range("E11").Value = '  +1.12%  '

# Row 12 - TRON
$ws.Range("E12").Value = '  -0.68%  '

# Row 13 - Toncoin
$ws.Range("E13").Value = '  +6.02%  '

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = '3.599.50'
$ws.Range("E14").Value = '  +3.78%  '

# Row 15 - Avalanche
$ws.Range("D15").Value = "'29.10"
$ws.Range("E15").Value = '  +3.77%  '

# Row 16 - WrappedBTC
$ws.Range("D16").Value = '76.272.26'
$ws.Range("E16").Value = '  +0.51%  '

# Row 17 - ShibaInu
$ws.Range("E17").Value = '  +1.89%  '

# Row 18 - WrappedEther
$ws.Range("D18").Value = '3.028.94'
$ws.Range("E18").Value = '  +3.21%  '

# Row 19 - Chainlink
$ws.Range("D19").Value = "'13.56"
$ws.Range("E19").Value = '  +2.72%  '

# Row 20 - Uniswap
$ws.Range("D20").Value = "'8.97"
$ws.Range("E20").Value = '  +1.14%  '

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'374.84"
$ws.Range("E21").Value = '  +0.70%  '

# Row 22 - SuiNetwork
$ws.Range("D22").Value = "'2.31"
$ws.Range("E22").Value = '  +0.40%  '

# Row 23 - Polkadot
$ws.Range("E23").Value = '  +1.45%  '

# Row 24 - Litecoin
$ws.Range("D24").Value = "'73.16"
$ws.Range("E24").Value = '  +2.07%  '

# Row 25 - WrappedeETH
$ws.Range("E25").Value = '  +3.68%  '

# Row 26 - Dai
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = '  +0.00%  '

# Row 27 - NEARProtocol
$ws.Range("D27").Value = "'4.36"
$ws.Range("E27").Value = '  +1.78%  '

# Row 28 - Aptos
$ws.Range("E28").Value = '  +1.50%  '

# Row 29 - PEPE
$ws.Range("E29").Value = '  +0.43%  '

# Row 30 - Binance-PegBSC-USD
$ws.Range("E30").Value = '  -0.05%  '

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "'8.29"
$ws.Range("E31").Value = '  +6.57%  '

# Row 32 - Fetch.AI
$ws.Range("E32").Value = '  +0.77%  '

# Row 33 - PancakeSwap
$ws.Range("D33").Value = "'1.94"
$ws.Range("E33").Value = '  +5.73%  '

# Row 34 - Bittensor
$ws.Range("D34").Value = "'491.53"
$ws.Range("E34").Value = '  -2.13%  '

# Row 35 - FirstDigitalUSD
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = '  +0.03%  '

# Row 36 - EthereumClassic
$ws.Range("D36").Value = "'20.66"
$ws.Range("E36").Value = '  +1.95%  '

# Row 37 - Monero
$ws.Range("D37").Value = "'162.71"
$ws.Range("E37").Value = '  -0.73%  '

# Row 38 - WhiteBITCoin
$ws.Range("E38").Value = '  +2.11%  '

# Row 39 - PolygonEcosystemToken
$ws.Range("D39").Value = "'0.384"
$ws.Range("E39").Value = '  +2.66%  '

# Row 40 - was Aave, now Kaspa (rank swap with row 41)
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = "'0.115"
$ws.Range("E40").Value = '  +1.99%  '

# Row 41 - was Kaspa, now Aave (rank swap with row 40)
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").Value = "'189.95"
$ws.Range("E41").Value = '  +4.81%  '

# Row 42 - Cronos
$ws.Range("D42").Value = "'0.104"
$ws.Range("E42").Value = '  -6.91%  '

# Row 43 - USDe
$ws.Range("E43").Value = '  -0.01%  '

# Row 44 - Mantle
$ws.Range("D44").Value = "'0.801"
$ws.Range("E44").Value = '  +21.92%  '

# Row 45 - RenderToken
$ws.Range("D45").Value = "'5.13"
$ws.Range("E45").Value = '  +2.65%  '

# Row 46 - ImmutableX
$ws.Range("E46").Value = '  +5.33%  '

# Row 47 - OKB
$ws.Range("D47").Value = "'42.05"
$ws.Range("E47").Value = '  +4.79%  '

# Row 48 - Stacks
$ws.Range("E48").Value = '  -0.60%  '

# Row 49 - dogwifhat
$ws.Range("E49").Value = '  +5.02%  '

# Row 50 - ARBITRUM
$ws.Range("E50").Value = '  +4.14%  '

# Row 51 - Filecoin
$ws.Range("E51").Value = '  +4.29%  '
